# Add a new "Date Created (Year)*" column (E) to Sheet1 with sample data,
# matching the dev-branch merge that introduced the image manifest's
# empty/extra column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "Date Created (Year)*"
$ws.Range("E2").Value = 2000
$ws.Range("E3").Value = 2000

# Style the new column's header + values (and the trailing blank row) with
# an explicit black font color, as applied upstream.
$ws.Range("E1:E4").Font.Color = 0

# Move the active selection to D1 (matches the merged sheet view state).
$ws.Range("D1").Select()
